# Commit: "Fruta / hortaliza, semanal"
# This weekly update inserts one new price-report row (row 90) into the
# "Arveja Verde" sheet, pushing the existing rows 90-164 down to 91-165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 90 (shifts rows 90:164 down to 91:165)
$ws.Rows("90:90").Insert()

# Populate the newly inserted row 90 with the new weekly data point
$ws.Range("A90").Value = 9
$ws.Range("B90").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C90").Value = "Metropolitana"
$ws.Range("D90").Value = 45090
$ws.Range("E90").Value = 13
$ws.Range("F90").Value = 100112022
$ws.Range("G90").Value = "Arveja Verde"
$ws.Range("H90").Value = "Perfection"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 43
$ws.Range("K90").Value = 36000
$ws.Range("L90").Value = 38000
$ws.Range("M90").Value = 37023
$ws.Range("N90").Value = "`$/saco 25 kilos"
$ws.Range("O90").Value = "Provincia de Huasco"
$ws.Range("P90").Value = 1481
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = "Hortaliza"
